# "Fruta / hortaliza, semanal"
# The data rows (2..34) are re-sorted (weekly ordering). Only the
# per-record measurement columns move between rows; every other column
# (A,B,C,E,F,G,H,I,N,O,Q,R) already matches across the whole block, so the
# edit is expressed as a row permutation over D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 34
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot the current ("before") values for the moving columns, keyed by row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# New row order: target row -> source row (1-indexed, matches the sheet).
$rowMap = @{
    2 = 12; 3 = 6; 4 = 25; 5 = 17; 6 = 33; 7 = 11; 8 = 9; 9 = 24; 10 = 31;
    11 = 26; 12 = 2; 13 = 7; 14 = 5; 15 = 8; 16 = 19; 17 = 15; 18 = 10;
    19 = 4; 20 = 22; 21 = 29; 22 = 27; 23 = 28; 24 = 3; 25 = 21; 26 = 14;
    27 = 18; 28 = 13; 29 = 30; 30 = 32; 31 = 23; 32 = 34; 33 = 16; 34 = 20
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $sourceVals[$c]
    }
}
